# Finished Week 13 logging
# Update the Target Depth Data for the "R" (Road?) row on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 426
$wsOff.Range("C3").Value = 300
$wsOff.Range("D3").Value = 97
$wsOff.Range("E3").Value = 54

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 391
$wsDef.Range("C3").Value = 272
$wsDef.Range("D3").Value = 105
$wsDef.Range("E3").Value = 52
$wsDef.Range("F3").Value = 9
